# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp on A1
# - Re-rank several countries (their name swaps to a different row while
#   the row that used to hold that name now shows another country's name)
# - Refresh various case-count figures (columns B-H) across many rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 09:05"

# --- Country re-ranking (column A) ------------------------------------
# El Salvador moves up in front of Eslovaquia; Eslovaquia/Eslovenia/
# Somalia/Gabon each shift one row down (rows 94-98).
$ws.Range("A94").Value = "El Salvador"
$ws.Range("A95").Value = "Eslovaquia"
$ws.Range("A96").Value = "Eslovenia"
$ws.Range("A97").Value = "Somalia"
$ws.Range("A98").Value = "Gabon"

# Nueva Caledonia and Belice swap places (rows 196-197).
$ws.Range("A196").Value = "Nueva Caledonia"
$ws.Range("A197").Value = "Belice"

# Seychelles moves up in front of Groenlandia; Groenlandia/Montserrat
# each shift one row down (rows 209-211).
$ws.Range("A209").Value = "Seychelles"
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Montserrat"

# San Bartolome and Bonaire, San Eustaquio y Saba swap places (rows 215-216).
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"

# --- Updated case figures (columns B-H) -------------------------------
$ws.Range("D11").Value = 155700
$ws.Range("E11").Value = 13466

$ws.Range("B94").Value = 1498
$ws.Range("C94").Value = 85
$ws.Range("D94").Value = 502
$ws.Range("E94").Value = 966
$ws.Range("H94").Value = 30

$ws.Range("B95").Value = 1495
$ws.Range("D95").Value = 1185
$ws.Range("E95").Value = 282
$ws.Range("H95").Value = 28

$ws.Range("B96").Value = 1466
$ws.Range("D96").Value = 1335
$ws.Range("E96").Value = 27
$ws.Range("H96").Value = 104

$ws.Range("B97").Value = 1455
$ws.Range("D97").Value = 163
$ws.Range("E97").Value = 1235
$ws.Range("H97").Value = 57

$ws.Range("B98").Value = 1432
$ws.Range("D98").Value = 301
$ws.Range("E98").Value = 1120
$ws.Range("H98").Value = 11

$ws.Range("D101").Value = 80
$ws.Range("E101").Value = 1022

$ws.Range("B121").Value = 702
$ws.Range("C121").Value = 1
$ws.Range("D121").Value = 456
$ws.Range("E121").Value = 234

$ws.Range("D132").Value = 401
$ws.Range("E132").Value = 32

$ws.Range("D196").Value = 18
$ws.Range("H196").Value = 0

$ws.Range("D197").Value = 16
$ws.Range("H197").Value = 2

$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1
